# Apply odds updates for Jogos_da_Semana_FlashScore_2024-10-13.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.7
$ws.Range("I2").Value = 5.75
$ws.Range("J2").Value = 2.4
$ws.Range("K2").Value = 1.92
$ws.Range("L2").Value = 6.5
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.5
$ws.Range("Z2").Value = 12
$ws.Range("AA2").Value = 17
$ws.Range("AG2").Value = 11
$ws.Range("AI2").Value = 21
$ws.Range("AN2").Value = 3.4
$ws.Range("AS2").Value = 251
$ws.Range("G5").Value = 2
$ws.Range("S5").Value = 1.58
$ws.Range("G7").Value = 1.33
$ws.Range("H7").Value = 4.7
$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 1.78
$ws.Range("K7").Value = 2.45
$ws.Range("L7").Value = 6.7
$ws.Range("N7").Value = 13.9
$ws.Range("Q7").Value = 1.53
$ws.Range("R7").Value = 2.2
$ws.Range("U7").Value = 1.8
$ws.Range("V7").Value = 1.8
$ws.Range("X7").Value = 6.8
$ws.Range("Z7").Value = 8.75
$ws.Range("AA7").Value = 10.5
$ws.Range("AD7").Value = 9.5
$ws.Range("AE7").Value = 19
$ws.Range("AG7").Value = 25
$ws.Range("AH7").Value = 60
$ws.Range("AI7").Value = 25
$ws.Range("AJ7").Value = 200
$ws.Range("AK7").Value = 90
$ws.Range("AL7").Value = 65
$ws.Range("AM7").Value = 500
$ws.Range("AN7").Value = 3.2
$ws.Range("AO7").Value = 5.9
$ws.Range("AQ7").Value = 15.5
$ws.Range("AU7").Value = 8
$ws.Range("AW7").Value = 8.75
$ws.Range("AX7").Value = 45
$ws.Range("AY7").Value = 40
$ws.Range("AZ7").Value = 300
$ws.Range("S10").Value = 1.41
$ws.Range("T10").Value = 2.62
$ws.Range("G23").Value = 2.7
$ws.Range("H23").Value = 3.8
$ws.Range("I23").Value = 2.25
$ws.Range("J23").Value = 3.2
$ws.Range("L23").Value = 2.88
$ws.Range("Q23").Value = 1.5
$ws.Range("X23").Value = 17
$ws.Range("Y23").Value = 11
$ws.Range("Z23").Value = 29
$ws.Range("AB23").Value = 23
$ws.Range("AI23").Value = 9.5
$ws.Range("AJ23").Value = 23
$ws.Range("AO23").Value = 15
$ws.Range("AW23").Value = 4.75
$ws.Range("AX23").Value = 12
$ws.Range("AY23").Value = 17
$ws.Range("BD23").Value = 151
$ws.Range("H24").Value = 3.75
$ws.Range("I24").Value = 4
$ws.Range("J24").Value = 2.3
$ws.Range("K24").Value = 2.22
$ws.Range("P24").Value = 3.7
$ws.Range("Q24").Value = 1.72
$ws.Range("R24").Value = 2.05
$ws.Range("X24").Value = 8.75
$ws.Range("AB24").Value = 23
$ws.Range("AD24").Value = 7.5
$ws.Range("AE24").Value = 14.5
$ws.Range("AF24").Value = 60
$ws.Range("AG24").Value = 13
$ws.Range("AH24").Value = 24
$ws.Range("AK24").Value = 37
$ws.Range("AL24").Value = 40
$ws.Range("AN24").Value = 3.65
$ws.Range("AO24").Value = 8.5
$ws.Range("AS24").Value = 250
$ws.Range("AU24").Value = 7.4
$ws.Range("AY24").Value = 29
$ws.Range("BA24").Value = 175
$ws.Range("BB24").Value = 400
$ws.Range("J26").Value = 3.6
$ws.Range("L26").Value = 3.3
$ws.Range("M26").Value = 1.06
$ws.Range("N26").Value = 5.8
$ws.Range("S26").Value = 1.52
$ws.Range("T26").Value = 2.22
$ws.Range("U26").Value = 1.93
$ws.Range("W26").Value = 6.7
$ws.Range("X26").Value = 13.5
$ws.Range("Y26").Value = 11
$ws.Range("AA26").Value = 32
$ws.Range("AD26").Value = 5.4
$ws.Range("AF26").Value = 100
$ws.Range("AG26").Value = 6.8
$ws.Range("AI26").Value = 10
$ws.Range("AJ26").Value = 32
$ws.Range("AK26").Value = 27
$ws.Range("AL26").Value = 40
$ws.Range("AO26").Value = 17
$ws.Range("AP26").Value = 27
$ws.Range("AQ26").Value = 90
$ws.Range("AR26").Value = 150
$ws.Range("AS26").Value = 400
$ws.Range("AT26").Value = 2.2
$ws.Range("AX26").Value = 15
$ws.Range("AY26").Value = 24
$ws.Range("AZ26").Value = 70
$ws.Range("BA26").Value = 110
$ws.Range("BB26").Value = 350
